$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4441.3335
$ws.Cells.Item(76, 9).Value = 4213.5
$ws.Cells.Item(76, 11).Value = 4213.5
$ws.Cells.Item(76, 13).Value = -3898.5

$ws.Cells.Item(79, 8).Value = 4441.3335
$ws.Cells.Item(79, 9).Value = 4213.5
$ws.Cells.Item(79, 11).Value = 4213.5
$ws.Cells.Item(79, 13).Value = -3121.5

$ws.Cells.Item(94, 8).Value = 167066670
$ws.Cells.Item(94, 9).Value = 500000000
$ws.Cells.Item(94, 11).Value = 500000000
$ws.Cells.Item(94, 13).Value = -499999549

$ws.Cells.Item(111, 8).Value = 1189.7142
$ws.Cells.Item(111, 9).Value = 1253.4
$ws.Cells.Item(111, 11).Value = 3760.2
$ws.Cells.Item(111, 13).Value = -693.2000000000003

$ws.Cells.Item(116, 8).Value = 20872088
$ws.Cells.Item(116, 9).Value = 20872088
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 20872088
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -20868646
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(138, 8).Value = 331778.88
$ws.Cells.Item(138, 9).Value = 627778.1
$ws.Cells.Item(138, 10).Value = 4621.8423
$ws.Cells.Item(138, 11).Value = 1883334.3
$ws.Cells.Item(138, 12).Value = 13865.5269
$ws.Cells.Item(138, 13).Value = -1878194.3
$ws.Cells.Item(138, 14).Value = -24145.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6750.7456
$ws.Cells.Item(32, 9).Value = 6635.185
$ws.Cells.Item(32, 10).Value = 7998.8
$ws.Cells.Item(32, 11).Value = 6635.185
$ws.Cells.Item(32, 12).Value = 7998.8
$ws.Cells.Item(32, 13).Value = -6348.185
$ws.Cells.Item(32, 14).Value = -8572.799999999999

$ws.Cells.Item(110, 8).Value = 2524.2354
$ws.Cells.Item(110, 9).Value = 1935
$ws.Cells.Item(110, 11).Value = 1935
$ws.Cells.Item(110, 13).Value = 110

$ws.Cells.Item(135, 8).Value = 34997.75
$ws.Cells.Item(135, 10).Value = 34997.75
$ws.Cells.Item(135, 12).Value = 34997.75
$ws.Cells.Item(135, 14).Value = -45137.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4817
$ws.Cells.Item(86, 9).Value = 8374.75
$ws.Cells.Item(86, 10).Value = 1767.5
$ws.Cells.Item(86, 11).Value = 8374.75
$ws.Cells.Item(86, 12).Value = 1767.5
$ws.Cells.Item(86, 13).Value = -7251.75
$ws.Cells.Item(86, 14).Value = -4013.5

$ws.Cells.Item(89, 8).Value = 4817
$ws.Cells.Item(89, 9).Value = 8374.75
$ws.Cells.Item(89, 10).Value = 1767.5
$ws.Cells.Item(89, 11).Value = 41873.75
$ws.Cells.Item(89, 12).Value = 8837.5
$ws.Cells.Item(89, 13).Value = -36257.75
$ws.Cells.Item(89, 14).Value = -20069.5

$ws.Cells.Item(105, 8).Value = 71023.734
$ws.Cells.Item(105, 9).Value = 112984
$ws.Cells.Item(105, 10).Value = 8083.3335
$ws.Cells.Item(105, 11).Value = 112984
$ws.Cells.Item(105, 12).Value = 8083.3335
$ws.Cells.Item(105, 13).Value = -111237
$ws.Cells.Item(105, 14).Value = -11577.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6870
$ws.Cells.Item(31, 9).Value = 7398.591
$ws.Cells.Item(31, 10).Value = 5812.8184
$ws.Cells.Item(31, 11).Value = 7398.591
$ws.Cells.Item(31, 12).Value = 5812.8184
$ws.Cells.Item(31, 13).Value = -7103.591
$ws.Cells.Item(31, 14).Value = -6402.8184

$ws.Cells.Item(34, 8).Value = 6870
$ws.Cells.Item(34, 9).Value = 7398.591
$ws.Cells.Item(34, 10).Value = 5812.8184
$ws.Cells.Item(34, 11).Value = 7398.591
$ws.Cells.Item(34, 12).Value = 5812.8184
$ws.Cells.Item(34, 13).Value = -7196.591
$ws.Cells.Item(34, 14).Value = -6216.8184

$ws.Cells.Item(97, 8).Value = 43500
$ws.Cells.Item(97, 10).Value = 43500
$ws.Cells.Item(97, 12).Value = 43500
$ws.Cells.Item(97, 14).Value = -45482

$ws.Cells.Item(122, 8).Value = 9310.799999999999
$ws.Cells.Item(122, 9).Value = 7873.9443
$ws.Cells.Item(122, 11).Value = 23621.8329
$ws.Cells.Item(122, 13).Value = -21171.8329

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 94864.39999999999
$ws.Cells.Item(37, 10).Value = 94864.39999999999
$ws.Cells.Item(37, 12).Value = 284593.2
$ws.Cells.Item(37, 14).Value = -284817.2

$ws.Cells.Item(103, 8).Value = 4950.2144
$ws.Cells.Item(103, 9).Value = 9026.333000000001
$ws.Cells.Item(103, 11).Value = 27078.999
$ws.Cells.Item(103, 13).Value = -26199.999

$ws.Cells.Item(137, 8).Value = 1736.75
$ws.Cells.Item(137, 9).Value = 1736.75
$ws.Cells.Item(137, 11).Value = 5210.25
$ws.Cells.Item(137, 13).Value = -110.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8784.875
$ws.Cells.Item(80, 9).Value = 10565.75
$ws.Cells.Item(80, 10).Value = 3442.25
$ws.Cells.Item(80, 11).Value = 10565.75
$ws.Cells.Item(80, 12).Value = 3442.25
$ws.Cells.Item(80, 13).Value = -9567.75
$ws.Cells.Item(80, 14).Value = -5438.25

$ws.Cells.Item(83, 8).Value = 8784.875
$ws.Cells.Item(83, 9).Value = 10565.75
$ws.Cells.Item(83, 10).Value = 3442.25
$ws.Cells.Item(83, 11).Value = 52828.75
$ws.Cells.Item(83, 12).Value = 17211.25
$ws.Cells.Item(83, 13).Value = -47836.75
$ws.Cells.Item(83, 14).Value = -27195.25

$ws.Cells.Item(102, 8).Value = 7094.9614
$ws.Cells.Item(102, 9).Value = 8493.223
$ws.Cells.Item(102, 11).Value = 8493.223
$ws.Cells.Item(102, 13).Value = -6871.223

$ws.Cells.Item(107, 8).Value = 322.86667
$ws.Cells.Item(107, 9).Value = 376.1111
$ws.Cells.Item(107, 10).Value = 243
$ws.Cells.Item(107, 11).Value = 376.1111
$ws.Cells.Item(107, 12).Value = 243
$ws.Cells.Item(107, 13).Value = 1543.8889
$ws.Cells.Item(107, 14).Value = -4083

$ws.Cells.Item(126, 8).Value = 8276.264999999999
$ws.Cells.Item(126, 9).Value = 9582.3125
$ws.Cells.Item(126, 10).Value = 7115.3335
$ws.Cells.Item(126, 11).Value = 28746.9375
$ws.Cells.Item(126, 12).Value = 21346.0005
$ws.Cells.Item(126, 13).Value = -26276.9375
$ws.Cells.Item(126, 14).Value = -26286.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4370.59
$ws.Cells.Item(122, 9).Value = 4216.5
$ws.Cells.Item(122, 11).Value = 12649.5
$ws.Cells.Item(122, 13).Value = -10199.5

$ws.Cells.Item(132, 8).Value = 712544.9
$ws.Cells.Item(132, 9).Value = 1066085.6
$ws.Cells.Item(132, 11).Value = 3198256.8
$ws.Cells.Item(132, 13).Value = -3195726.8

$ws.Cells.Item(135, 8).Value = 97997.164
$ws.Cells.Item(135, 10).Value = 97997.164
$ws.Cells.Item(135, 12).Value = 97997.164
$ws.Cells.Item(135, 14).Value = -108137.164

$ws.Cells.Item(140, 8).Value = 84249.836
$ws.Cells.Item(140, 10).Value = 99899.8
$ws.Cells.Item(140, 12).Value = 99899.8
$ws.Cells.Item(140, 14).Value = -110259.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 21551.4
$ws.Cells.Item(107, 9).Value = 1893.0834
$ws.Cells.Item(107, 10).Value = 100184.664
$ws.Cells.Item(107, 11).Value = 5679.2502
$ws.Cells.Item(107, 12).Value = 300553.992
$ws.Cells.Item(107, 13).Value = -3759.2502
$ws.Cells.Item(107, 14).Value = -304393.992

$ws.Cells.Item(135, 8).Value = 60000
$ws.Cells.Item(135, 10).Value = 60000
$ws.Cells.Item(135, 12).Value = 60000
$ws.Cells.Item(135, 14).Value = -70140

$ws.Cells.Item(136, 8).Value = 561347.4399999999
$ws.Cells.Item(136, 9).Value = 708510.4399999999
$ws.Cells.Item(136, 11).Value = 2125531.32
$ws.Cells.Item(136, 13).Value = -2122981.32
